# US3 and 4 done
# Fills in the "24. oktober" row (adds " og bug fix" + "10 timer") and the
# "25. oktober" row (date, description with a spell-checked "backend", and
# "12 timer") of the Tidskema table.

$d = $word.ActiveDocument
$table = $d.Tables.Item(1)

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# --- Row 5 ("24. oktober") ------------------------------------------------

# Description cell: append a new run " og bug fix" after the existing text,
# keeping "Udarbejdelse af frontend" as its own separate run.
$descCell = $table.Cell(5, 2)
$descXml = "<w:p $wNs>" +
  "<w:r><w:t>Udarbejdelse af frontend</w:t></w:r>" +
  "<w:r><w:t xml:space=`"preserve`"> og bug fix</w:t></w:r>" +
  "</w:p>"
$descCell.Range.InsertXML($descXml)

# Time spent cell: was empty, now "10 timer".
$table.Cell(5, 3).Range.Text = "10 timer"

# --- Row 6 ("25. oktober") -------------------------------------------------

# Date cell.
$table.Cell(6, 1).Range.Text = "25. oktober"

# Description cell: includes a spell-check bracketed "backend".
$descCell2 = $table.Cell(6, 2)
$descXml2 = "<w:p $wNs>" +
  "<w:r><w:t xml:space=`"preserve`">Fortsat udarbejdelse af frontend og rettelse af </w:t></w:r>" +
  "<w:proofErr w:type=`"spellStart`"/>" +
  "<w:r><w:t>backend</w:t></w:r>" +
  "<w:proofErr w:type=`"spellEnd`"/>" +
  "<w:r><w:t xml:space=`"preserve`"> fejl</w:t></w:r>" +
  "</w:p>"
$descCell2.Range.InsertXML($descXml2)

# Time spent cell.
$table.Cell(6, 3).Range.Text = "12 timer"
